# "Hjemme passive tweaks lichtwark deleted values"
# Update the B1:E1 header counts and the B2:E3 data block on Ark1, then
# restore the selection to the (now smaller) edited range B1:E3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (count of samples per column).
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 data.
$ws.Range("B2").Value = 257.38178829125468
$ws.Range("C2").Value = 302.84293275323353
$ws.Range("D2").Value = 253.51132356186886
$ws.Range("E2").Value = 304.56824208576404

# Row 3 data.
$ws.Range("B3").Value = 249.3079991300373
$ws.Range("C3").Value = 301.16146993977605
$ws.Range("D3").Value = 247.47747324598532
$ws.Range("E3").Value = 307.55551255631553

# Selection now only spans the edited block instead of the whole table.
$ws.Range("B1:E3").Select()
